# Applies the "Update the report and add stories in ged file" commit.
# - Adds two new rows of test-coverage data to the Sprint3 sheet
#   (user stories "Fewer than 15 siblings" and "Siblings should not marry").
# - Updates the Sprint3 sheet view (zoom / selection).
# - Refreshes row heights across several sheets (Backlog, Sprint1, Sprint2,
#   Sprint3, Stories) to reflect the re-saved layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sprint3: fill in newly completed test-coverage details for two user stories
# ---------------------------------------------------------------------------
$ws = $wb.Sheets.Item("Sprint3")

# Row 2 - US15 "Fewer than 15 siblings"
$ws.Range("G2").Value = 12
$ws.Range("H2").Value = 1
$ws.Range("I2").NumberFormat = "m/d/yyyy"
$ws.Range("I2").VerticalAlignment = -4160
$ws.Range("I2").WrapText = $true
$ws.Range("I2").Value = 42092
$ws.Range("J2").Value = "siblings_not_too_many.py"
$ws.Range("K2").Value = "siblings_not_too_many"
$ws.Range("L2").Value = 15
$ws.Range("N2").Value = "test_siblings_not_too_many.py"
$ws.Range("O2").Value = "test_valid_siblings_num, test_siblings_more_than_15"
$ws.Range("P2").Value = 42

# Row 5 - US18 "Siblings should not marry"
$ws.Range("G5").Value = 18
$ws.Range("H5").Value = 2
$ws.Range("I5").NumberFormat = "m/d/yyyy"
$ws.Range("I5").VerticalAlignment = -4160
$ws.Range("I5").WrapText = $true
$ws.Range("I5").Value = 42092
$ws.Range("J5").Value = "siblings_not_marry.py"
$ws.Range("K5").Value = "siblings_not_marry"
$ws.Range("L5").Value = 18
$ws.Range("N5").Value = "test_siblings_not_marry.py"
$ws.Range("O5").Value = "test_valid, test_invalid"
$ws.Range("P5").Value = 264

# Update the sheet view: zoom to 90%, select P5, and drop the old
# top-left-cell / selection anchored at P9.
[void]$ws.Activate()
$excel.ActiveWindow.Zoom = 90
$ws.Range("P5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Row-height refresh (auto-fit heights recorded after the re-save)
# ---------------------------------------------------------------------------
$ws = $wb.Sheets.Item("Backlog")
$ws.Range("A2:A4").EntireRow.RowHeight = 34
$ws.Rows.Item(5).RowHeight = 51
$ws.Range("A6:A7").EntireRow.RowHeight = 34
$ws.Rows.Item(8).RowHeight = 68
$ws.Range("A9:A11").EntireRow.RowHeight = 51
$ws.Rows.Item(12).RowHeight = 34
$ws.Rows.Item(13).RowHeight = 68
$ws.Rows.Item(14).RowHeight = 85
$ws.Range("A15:A22").EntireRow.RowHeight = 34
$ws.Rows.Item(23).RowHeight = 51
$ws.Rows.Item(24).RowHeight = 68
$ws.Rows.Item(25).RowHeight = 51
$ws.Range("A26:A32").EntireRow.RowHeight = 34
$ws.Rows.Item(33).RowHeight = 51

$ws = $wb.Sheets.Item("Sprint1")
$ws.Rows.Item(1).RowHeight = 14
$ws.Range("A2:A7").EntireRow.RowHeight = 28
$ws.Rows.Item(8).RowHeight = 42
$ws.Rows.Item(9).RowHeight = 28
$ws.Rows.Item(13).RowHeight = 14
$ws.Rows.Item(15).RowHeight = 14
$ws.Rows.Item(16).RowHeight = 24
$ws.Rows.Item(18).RowHeight = 14
$ws.Rows.Item(19).RowHeight = 72
$ws.Rows.Item(21).RowHeight = 14

$ws = $wb.Sheets.Item("Sprint2")
$ws.Rows.Item(2).RowHeight = 84
$ws.Rows.Item(3).RowHeight = 98
$ws.Rows.Item(4).RowHeight = 42
$ws.Range("A5:A6").EntireRow.RowHeight = 28
$ws.Rows.Item(7).RowHeight = 42
$ws.Range("A8:A9").EntireRow.RowHeight = 28
$ws.Rows.Item(14).RowHeight = 14
$ws.Rows.Item(16).RowHeight = 28
$ws.Rows.Item(17).RowHeight = 24
$ws.Rows.Item(19).RowHeight = 28
$ws.Rows.Item(20).RowHeight = 72
$ws.Rows.Item(22).RowHeight = 14
$ws.Rows.Item(23).RowHeight = 36

$ws = $wb.Sheets.Item("Sprint3")
$ws.Rows.Item(2).RowHeight = 84
$ws.Range("A3:A4").EntireRow.RowHeight = 14
$ws.Range("A5:A6").EntireRow.RowHeight = 42
$ws.Range("A7:A8").EntireRow.RowHeight = 224
$ws.Rows.Item(9).RowHeight = 42

$ws = $wb.Sheets.Item("Stories")
$ws.Rows.Item(2).RowHeight = 34
$ws.Range("A3:A4").EntireRow.RowHeight = 17
$ws.Rows.Item(5).RowHeight = 34
$ws.Range("A6:A7").EntireRow.RowHeight = 17
$ws.Rows.Item(8).RowHeight = 51
$ws.Range("A9:A12").EntireRow.RowHeight = 34
$ws.Rows.Item(13).RowHeight = 51
$ws.Rows.Item(14).RowHeight = 68
$ws.Rows.Item(15).RowHeight = 34
$ws.Rows.Item(16).RowHeight = 17
$ws.Rows.Item(17).RowHeight = 34
$ws.Range("A18:A20").EntireRow.RowHeight = 17
$ws.Range("A21:A24").EntireRow.RowHeight = 34
$ws.Rows.Item(25).RowHeight = 51
$ws.Rows.Item(26).RowHeight = 34
$ws.Rows.Item(27).RowHeight = 136
$ws.Range("A28:A29").EntireRow.RowHeight = 34
$ws.Range("A30:A31").EntireRow.RowHeight = 17
$ws.Rows.Item(32).RowHeight = 34
$ws.Rows.Item(33).RowHeight = 17
$ws.Rows.Item(34).RowHeight = 34
$ws.Rows.Item(35).RowHeight = 51
$ws.Range("A36:A43").EntireRow.RowHeight = 34
